$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset grew by one weekly observation: a new record is inserted at
# row 43, which pushes the former rows 43-53 down to rows 44-54 (and the
# sheet's used range grows from A1:T53 to A1:T54).
$ws.Rows(43).Insert()

$ws.Cells.Item(43, 1).Value = 6
$ws.Cells.Item(43, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(43, 3).Value = "Metropolitana"
$ws.Cells.Item(43, 4).Value = 44455
$ws.Cells.Item(43, 5).Value = 13
$ws.Cells.Item(43, 6).Value = "Fruta"
$ws.Cells.Item(43, 7).Value = 100108
$ws.Cells.Item(43, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(43, 9).Value = 100108007
$ws.Cells.Item(43, 10).Value = "Coco"
$ws.Cells.Item(43, 11).Value = "Sin especificar"
$ws.Cells.Item(43, 12).Value = "Primera"
$ws.Cells.Item(43, 13).Value = 250
$ws.Cells.Item(43, 14).Value = 20000
$ws.Cells.Item(43, 15).Value = 20000
$ws.Cells.Item(43, 16).Value = 20000
$ws.Cells.Item(43, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(43, 18).Value = "Perú"
$ws.Cells.Item(43, 19).Value = 1000
$ws.Cells.Item(43, 20).Value = 20
